# Insert a new week of "Sandia" price records (Primera/Segunda/Tercera,
# date 2023-01-17 / serial 44943) above the existing Región de O'Higgins
# block, pushing the later rows down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 132:148 down to 135:151, carrying their formatting along.
$ws.Rows("132:134").Insert()

# Fill in the shared / repeated columns for the three new rows.
$ws.Range("A132:A134").Value = 11
$ws.Range("B132:B134").Value = "Vega Monumental Concepción"
$ws.Range("C132:C134").Value = "Bíobío"
$ws.Range("D132:D134").Value = 44943
$ws.Range("E132:E134").Value = 8
$ws.Range("F132:F134").Value = 100112028
$ws.Range("G132:G134").Value = "Sandia"
$ws.Range("H132:H134").Value = "Sin especificar"
$ws.Range("N132:N134").Value = "$/unidad"
$ws.Range("O132:O134").Value = "Región de O'Higgins"
$ws.Range("Q132:Q134").Value = 1
$ws.Range("R132:R134").Value = "Hortaliza"

# Row 132: Primera
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2500
$ws.Range("L132").Value = 2500
$ws.Range("M132").Value = 2500
$ws.Range("P132").Value = 2500

# Row 133: Segunda
$ws.Range("I133").Value = "Segunda"
$ws.Range("J133").Value = 850
$ws.Range("K133").Value = 2000
$ws.Range("L133").Value = 2000
$ws.Range("M133").Value = 2000
$ws.Range("P133").Value = 2000

# Row 134: Tercera
$ws.Range("I134").Value = "Tercera"
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 1800
$ws.Range("L134").Value = 1800
$ws.Range("M134").Value = 1800
$ws.Range("P134").Value = 1800
